$d = $word.ActiveDocument

# Locate the paragraph that currently reads "ntre em Contato*" (the leading
# "E" is missing, which is why Word's proofer flagged "ntre" as a misspelled
# word with a spellStart/spellEnd pair splitting the run in two).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "ntre em Contato*") {
        # Replace the whole paragraph (runs only, keep the paragraph mark)
        # with the corrected text "Entre em Contato*", typed as two runs
        # ("E" then "ntre em Contato*") and with no leftover proofing marks,
        # matching how Word re-serializes the paragraph once the word is
        # spelled correctly again.
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>E</w:t></w:r><w:r><w:t>ntre em Contato*</w:t></w:r></w:p>'
        $p.Range.InsertXML($xml)
        break
    }
}
